# "Generate Report for Handoff" - refresh the localization-status report
# with a new source-file GUID, new xliff hashes and new handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "8b10ff1d-7258-479b-9cee-88ff74bfa152"
$newGuid = "32bc3f0e-38ea-4c03-9f6b-ed86db80e5dc"
$newXlfHash = "e98a497135d58d9a27ebdc4d713715c9775b040e"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aeca6f86891acf83507447d735b19c551748a951/e2e"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-16 18:54:34"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$repoBase/$newGuid.md", "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Style = "Normal"

$wsZh.Range("G2").Value = "$newGuid.$newXlfHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-16 18:54:29"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoBase/$newGuid.md", "", "", "$newGuid.md")

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Style = "Normal"

$wsDe.Range("G2").Value = "$newGuid.$newXlfHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-16 18:54:34"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBase/$newGuid.md", "", "", "$newGuid.md")

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
